$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 283; everything from the old
# row 283 downward (old rows 283-298) shifts down to become rows 285-300.
$ws.Rows.Item(283).Resize(2).Insert()

# New row 283: a new weekly price observation (docena de paquetes, Provincia de Cautín)
$ws.Cells.Item(283, 1).Value = 10
$ws.Cells.Item(283, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(283, 3).Value = "La Araucanía"
$ws.Cells.Item(283, 4).Value = 44516
$ws.Cells.Item(283, 5).Value = 9
$ws.Cells.Item(283, 6).Value = 100114014
$ws.Cells.Item(283, 7).Value = "Betarraga"
$ws.Cells.Item(283, 8).Value = "Sin especificar"
$ws.Cells.Item(283, 9).Value = "Primera"
$ws.Cells.Item(283, 10).Value = 120
$ws.Cells.Item(283, 11).Value = 9000
$ws.Cells.Item(283, 12).Value = 10000
$ws.Cells.Item(283, 13).Value = 9542
$ws.Cells.Item(283, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(283, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(283, 16).Value = 795
$ws.Cells.Item(283, 17).Value = 12
$ws.Cells.Item(283, 18).Value = "Hortaliza"

# New row 284: a new weekly price observation (paquete 5 unidades, Región Metropolitana)
$ws.Cells.Item(284, 1).Value = 10
$ws.Cells.Item(284, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(284, 3).Value = "La Araucanía"
$ws.Cells.Item(284, 4).Value = 44516
$ws.Cells.Item(284, 5).Value = 9
$ws.Cells.Item(284, 6).Value = 100114014
$ws.Cells.Item(284, 7).Value = "Betarraga"
$ws.Cells.Item(284, 8).Value = "Sin especificar"
$ws.Cells.Item(284, 9).Value = "Primera"
$ws.Cells.Item(284, 10).Value = 65
$ws.Cells.Item(284, 11).Value = 700
$ws.Cells.Item(284, 12).Value = 700
$ws.Cells.Item(284, 13).Value = 700
$ws.Cells.Item(284, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(284, 15).Value = "Región Metropolitana"
$ws.Cells.Item(284, 16).Value = 140
$ws.Cells.Item(284, 17).Value = 5
$ws.Cells.Item(284, 18).Value = "Hortaliza"
